# chitieupt-samplesheet.xlsx edit:
# Integrate "Phuong Phap" (method) as a field of "chitieuphantich" (indicator)
# placed between "Ten Chi Tieu" and "Don Gia" (swap columns C and D
# semantics), matching the new "Bang Bao Gia 2016" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): swap "Don Gia" / "Phuong Phap" ---
$ws.Range("C2").Value = "Phương Pháp"
$ws.Range("D2").Value = "Đơn Giá"

# --- Row 3 (Nuoc / abc): price moves C->D, method moves D->C ---
$ws.Range("C3").Value = "PP1"
$ws.Range("D3").Value = 100

# --- Row 4: extra method value, column D -> C ---
$ws.Range("C4").Value = "PP2"
$ws.Range("D4").ClearContents()

# --- Row 5: extra method value, column D -> C ---
$ws.Range("C5").Value = "PP3"
$ws.Range("D5").ClearContents()

# --- Row 6 (acc): price moves C->D, method moves D->C ---
$ws.Range("C6").Value = "PPACC1"
$ws.Range("D6").Value = 10000

# --- Row 7: extra method value, column D -> C ---
$ws.Range("C7").Value = "PPACC2"
$ws.Range("D7").ClearContents()

# --- Row 8 (Bun Thai / ccc): price moves C->D, method moves D->C ---
$ws.Range("C8").Value = "BTPP0"
$ws.Range("D8").Value = 1000000

# --- Row 9: extra method value, column D -> C ---
$ws.Range("C9").Value = "BTPP11"
$ws.Range("D9").ClearContents()

# --- Row 10 (sde): price moves C->D, method moves D->C ---
$ws.Range("C10").Value = "SDE11EEPP"
$ws.Range("D10").Value = 30000

# --- Row 11 (Khong Khi / fde): price moves C->D, method moves D->C ---
$ws.Range("C11").Value = "Lọc KK1"
$ws.Range("D11").Value = 20001

# --- Row 12: extra method value, column D -> C ---
$ws.Range("C12").Value = "Loc KKK"
$ws.Range("D12").ClearContents()

# --- Column widths: column C (Phuong Phap) and D (Don Gia) get their own
#     explicit widths now instead of D alone being wide ---
#     (target OOXML widths are 15.84 / 12.22 "characters"; the engine
#     quantizes ColumnWidth to 1/6-character pixel steps the same way
#     Excel itself does, so these are the closest achievable inputs)
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 11.333333333333334

# --- Font charset for the title font (Times New Roman) ---
$ws.Range("A1").Font.Charset = 1

# --- Selection moved from F8 to H10 ---
[void]$ws.Range("H10").Select()
